$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.167.40"
$ws.Range("E2").Value = "  -2.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.340.88"
$ws.Range("E3").Value = "  -4.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.84"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.69"
$ws.Range("E6").Value = "  -6.47%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").Value = "  -5.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.335.48"
$ws.Range("E9").Value = "  -4.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.72"
$ws.Range("E11").Value = "  -7.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.79"
$ws.Range("E13").Value = "  -9.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000258"
$ws.Range("E14").Value = "  -6.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.876.80"
$ws.Range("E15").Value = "  -4.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.74"
$ws.Range("E16").Value = "  -6.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.125.23"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.325.88"
$ws.Range("E18").Value = "  -5.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "561.51"
$ws.Range("E19").Value = "  -7.91%  "
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.20"
$ws.Range("E21").Value = "  -6.46%  "
$ws.Range("E22").Value = "  -7.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.50"
$ws.Range("E23").Value = "  -6.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "90.40"
$ws.Range("E24").Value = "  -7.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.29"
$ws.Range("E25").Value = "  -8.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.48"
$ws.Range("E26").Value = "  -6.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  -8.01%  "
$ws.Range("E29").Value = "  -10.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.20"
$ws.Range("E30").Value = "  -8.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.34"
$ws.Range("E31").Value = "  -9.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  -7.07%  "
$ws.Range("E33").Value = "  -10.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "580.62"
$ws.Range("E34").Value = "  -9.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.25"
$ws.Range("E35").Value = "  -9.13%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.00"
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0452"
$ws.Range("E38").Value = "  -3.83%  "
$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.94"
$ws.Range("E39").Value = "  -7.04%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0908"
$ws.Range("E40").Value = "  -8.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.137"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.96"
$ws.Range("E42").Value = "  -17.27%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.082.46"
$ws.Range("E43").Value = "  -8.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.66"
$ws.Range("E44").Value = "  -7.86%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0628"
$ws.Range("E45").Value = "  -15.45%  "
$ws.Range("E46").Value = "  -9.92%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.96"
$ws.Range("E47").Value = "  -9.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.26"
$ws.Range("E48").Value = "  -11.38%  "
$ws.Range("E49").Value = "  -6.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "129.95"
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("E51").Value = "  -0.03%  "
